$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.586692571640015
$ws.Range("B1").Value = 4.872288703918457
$ws.Range("C1").Value = 6.49193811416626
$ws.Range("D1").Value = 6.453382015228271
$ws.Range("E1").Value = 5.352084636688232
